$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 20002080
$ws.Range("B2").Value = "Phạm Hồng Nghĩa"
$ws.Range("C2").Value = 37536
$ws.Range("G2").Value = 3.77

$ws.Range("A3").Value = 20002076
$ws.Range("B3").Value = "Dương Văn Nam"
$ws.Range("C3").Value = 37536
$ws.Range("G3").Value = 3.73

$ws.Range("A4").Value = 20002077
$ws.Range("B4").Value = "Lã Đức Nam"
$ws.Range("C4").Value = 37588
$ws.Range("G4").Value = 3.65
